# Update "想去人数" (want-to-go count) figures on the 南宁-漫展信息 workbook.
# Affects the "展览" sheet and the combined "全部类型" sheet (F column values).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 13865
$ws1.Range("F3").Value = 327
$ws1.Range("F4").Value = 669
$ws1.Range("F5").Value = 235
$ws1.Range("F6").Value = 507
$ws1.Range("F7").Value = 1432
$ws1.Range("F8").Value = 136

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 13865
$ws4.Range("F3").Value = 327
$ws4.Range("F4").Value = 669
$ws4.Range("F5").Value = 235
$ws4.Range("F8").Value = 507
$ws4.Range("F9").Value = 1432
$ws4.Range("F11").Value = 136
